$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.667.82"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.279.83"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "95.76"
$ws.Range("E5").Value = "  -2.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "267.15"
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.609"
$ws.Range("E9").Value = "  -4.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.22"
$ws.Range("E10").Value = "  -7.93%  "
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.73"
$ws.Range("E12").Value = "  -7.83%  "
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.620.78"
$ws.Range("E14").Value = "  +1.68%  "
$ws.Range("E15").Value = "  -1.37%  "
$ws.Range("E16").Value = "  +2.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.284.85"
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.577.33"
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000108"
$ws.Range("E19").Value = "  +2.10%  "
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.94"
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.35"
$ws.Range("E22").Value = "  +1.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.79"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  -2.71%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.34"
$ws.Range("E26").Value = "  -1.55%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.49"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.44"
$ws.Range("E28").Value = "  -2.89%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.39"
$ws.Range("E30").Value = "  -4.94%  "
$ws.Range("E31").Value = "  +1.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.01"
$ws.Range("E32").Value = "  +4.48%  "
$ws.Range("E33").Value = "  -4.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.37"
$ws.Range("E34").Value = "  -3.94%  "
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("E36").Value = "  -3.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0356"
$ws.Range("E37").Value = "  +1.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.44"
$ws.Range("E38").Value = "  +2.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.28"
$ws.Range("E39").Value = "  -11.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.36"
$ws.Range("E40").Value = "  +7.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.238"
$ws.Range("E41").Value = "  -4.34%  "
$ws.Range("E42").Value = "  +17.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.89"
$ws.Range("E43").Value = "  -6.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.92"
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.83"
$ws.Range("E45").Value = "  +4.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.21"
$ws.Range("E46").Value = "  -4.08%  "
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.30"
$ws.Range("E48").Value = "  -2.24%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("E50").Value = "  +8.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.501.29"
$ws.Range("E51").Value = "  +1.59%  "
Write-Output "applied 85 changes"
